$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for 2022-Q3 at the top
#    of the data (row 2), pushing existing quarters down by one row,
#    and append the tail-end "2020-Q4" row that falls off the bottom.
# -----------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Row 2 lost its formatting on insert (Excel copies the header's bold/
# bordered style into the new row) -- pull the plain data-row styling
# back from row 3 (the old row 2, now shifted down) so column A keeps
# its bordered/bold look and B:D stay plain, matching the other rows.
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))
$summary.Cells.Item(3,2).Copy($summary.Cells.Item(2,2))
$summary.Cells.Item(3,3).Copy($summary.Cells.Item(2,3))
$summary.Cells.Item(3,4).Copy($summary.Cells.Item(2,4))

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 12
$summary.Cells.Item(2,4).Value = 1.02

# Column A is just a running 0-based index -- renumber all data rows
# now that an extra one exists.
for ($i = 0; $i -lt 8; $i++) {
    $summary.Cells.Item(2 + $i, 1).Value = $i
}

# -----------------------------------------------------------------
# 2) Add the new "2022-Q3" sheet (a fund-holdings breakdown), placed
#    right after "总计" and before "2022-Q2". Cloning "2022-Q2" gives
#    us the right header row + column styling for free; we then wipe
#    its sample data and refill with the 2022-Q3 numbers.
# -----------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

$q3.Range("A2:H7").ClearContents()

# Extend column A's bordered/bold style down through row 13 (12 data
# rows) and force columns B and D:G to text so fund codes keep their
# leading zeros and the percentage/NAV figures keep trailing zeros.
$q3.Cells.Item(2,1).Copy($q3.Range("A2:A13"))
$q3.Range("B2:B13").NumberFormat = "@"
$q3.Range("D2:G13").NumberFormat = "@"

$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "000979"
$q3.Cells.Item(2,3).Value = "景顺长城沪港深精选股票"
$q3.Cells.Item(2,4).Value = "20.32"
$q3.Cells.Item(2,5).Value = "80.04"
$q3.Cells.Item(2,6).Value = "1.98"
$q3.Cells.Item(2,7).Value = "0.4023"
$q3.Cells.Item(2,8).Value = 9
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "008850"
$q3.Cells.Item(3,3).Value = "景顺长城价值稳进三年定期开放灵活配置混合"
$q3.Cells.Item(3,4).Value = "17.29"
$q3.Cells.Item(3,5).Value = "84.84"
$q3.Cells.Item(3,6).Value = "1.98"
$q3.Cells.Item(3,7).Value = "0.3423"
$q3.Cells.Item(3,8).Value = 9
$q3.Cells.Item(4,1).Value = 2
$q3.Cells.Item(4,2).Value = "008060"
$q3.Cells.Item(4,3).Value = "景顺长城价值边际灵活配置混合A"
$q3.Cells.Item(4,4).Value = "5.45"
$q3.Cells.Item(4,5).Value = "80.34"
$q3.Cells.Item(4,6).Value = "2.02"
$q3.Cells.Item(4,7).Value = "0.1101"
$q3.Cells.Item(4,8).Value = 9
$q3.Cells.Item(5,1).Value = 3
$q3.Cells.Item(5,2).Value = "011583"
$q3.Cells.Item(5,3).Value = "大成港股精选混合（QDII）A"
$q3.Cells.Item(5,4).Value = "2.62"
$q3.Cells.Item(5,5).Value = "82.40"
$q3.Cells.Item(5,6).Value = "3.24"
$q3.Cells.Item(5,7).Value = "0.0849"
$q3.Cells.Item(5,8).Value = 9
$q3.Cells.Item(6,1).Value = 4
$q3.Cells.Item(6,2).Value = "501310"
$q3.Cells.Item(6,3).Value = "华宝标普沪港深中国增强价值指数（LOF）A"
$q3.Cells.Item(6,4).Value = "0.89"
$q3.Cells.Item(6,5).Value = "93.50"
$q3.Cells.Item(6,6).Value = "2.95"
$q3.Cells.Item(6,7).Value = "0.0263"
$q3.Cells.Item(6,8).Value = 8
$q3.Cells.Item(7,1).Value = 5
$q3.Cells.Item(7,2).Value = "015779"
$q3.Cells.Item(7,3).Value = "景顺长城价值边际灵活配置混合C"
$q3.Cells.Item(7,4).Value = "0.79"
$q3.Cells.Item(7,5).Value = "80.34"
$q3.Cells.Item(7,6).Value = "2.02"
$q3.Cells.Item(7,7).Value = "0.0160"
$q3.Cells.Item(7,8).Value = 9
$q3.Cells.Item(8,1).Value = 6
$q3.Cells.Item(8,2).Value = "161620"
$q3.Cells.Item(8,3).Value = "融通核心价值混合（QDII）A"
$q3.Cells.Item(8,4).Value = "0.55"
$q3.Cells.Item(8,5).Value = "57.96"
$q3.Cells.Item(8,6).Value = "2.73"
$q3.Cells.Item(8,7).Value = "0.0150"
$q3.Cells.Item(8,8).Value = 9
$q3.Cells.Item(9,1).Value = 7
$q3.Cells.Item(9,2).Value = "011584"
$q3.Cells.Item(9,3).Value = "大成港股精选混合（QDII）C"
$q3.Cells.Item(9,4).Value = "0.42"
$q3.Cells.Item(9,5).Value = "82.40"
$q3.Cells.Item(9,6).Value = "3.24"
$q3.Cells.Item(9,7).Value = "0.0136"
$q3.Cells.Item(9,8).Value = 9
$q3.Cells.Item(10,1).Value = 8
$q3.Cells.Item(10,2).Value = "011647"
$q3.Cells.Item(10,3).Value = "博时港股通红利精选混合A"
$q3.Cells.Item(10,4).Value = "0.11"
$q3.Cells.Item(10,5).Value = "82.44"
$q3.Cells.Item(10,6).Value = "4.58"
$q3.Cells.Item(10,7).Value = "0.0050"
$q3.Cells.Item(10,8).Value = 6
$q3.Cells.Item(11,1).Value = 9
$q3.Cells.Item(11,2).Value = "007397"
$q3.Cells.Item(11,3).Value = "华宝标普沪港深中国增强价值指数（LOF）C"
$q3.Cells.Item(11,4).Value = "0.04"
$q3.Cells.Item(11,5).Value = "93.50"
$q3.Cells.Item(11,6).Value = "2.95"
$q3.Cells.Item(11,7).Value = "0.0012"
$q3.Cells.Item(11,8).Value = 8
$q3.Cells.Item(12,1).Value = 10
$q3.Cells.Item(12,2).Value = "011648"
$q3.Cells.Item(12,3).Value = "博时港股通红利精选混合C"
$q3.Cells.Item(12,4).Value = "0.02"
$q3.Cells.Item(12,5).Value = "82.44"
$q3.Cells.Item(12,6).Value = "4.58"
$q3.Cells.Item(12,7).Value = "0.0009"
$q3.Cells.Item(12,8).Value = 6
$q3.Cells.Item(13,1).Value = 11
$q3.Cells.Item(13,2).Value = "014127"
$q3.Cells.Item(13,3).Value = "融通核心价值混合（QDII）C"
$q3.Cells.Item(13,4).Value = "0.01"
$q3.Cells.Item(13,5).Value = "57.96"
$q3.Cells.Item(13,6).Value = "2.73"
$q3.Cells.Item(13,7).Value = "0.0003"
$q3.Cells.Item(13,8).Value = 9

$summary.Select()
